# Rename the "_old"/"_new" suffixed column headers to "_FV2410"/"_FV2504"
# respectively (the "diff" header in the middle column stays unchanged),
# freeze the header row, and turn the used range into a native Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J (1-10): "<field>_old"  -> "<field>_FV2410"
for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = "$($fields[$i])_FV2410"
}

# Column K (11) stays "diff"

# Columns L-U (12-21): "<field>_new" -> "<field>_FV2504"
for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = "$($fields[$i])_FV2504"
}

# Freeze the header row (pane split after row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Convert the data range into an Excel Table ("Table1") with the renamed
# headers as its column names.
$tableRange = $ws.Range("A1:U77")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

Write-Output "done"
